# Generate Report for Archive
# Update localization status from "Ready for handoff" to "In Translation"
# across the Overview, zh-cn and de-de sheets, and shrink the now-narrower
# "Status" columns to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Narrow columns E and F to match the shorter replacement text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status cell (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status cell (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
